$wb = $excel.ActiveWorkbook

# --- Update localization status text: "Ready for handoff" -> "In Translation" ---
# This shared string is referenced from the Overview sheet (E2, F2) and from
# the per-language detail sheets (Status column, C2). Writing the same new
# text to every referencing cell collapses them back onto a single shared
# string entry, matching the original single shared-string edit.
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E2").Value = "In Translation"
$wsOverview.Range("F2").Value = "In Translation"

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C2").Value = "In Translation"

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C2").Value = "In Translation"

# --- Narrow the language/status columns ---
# Overview: zh-cn (E) and de-de (F) columns
$wsOverview.Columns.Item(5).ColumnWidth = 12.5
$wsOverview.Columns.Item(6).ColumnWidth = 12.5

# zh-cn / de-de detail sheets: Status column (C)
$wsZhCn.Columns.Item(3).ColumnWidth = 12.5
$wsDeDe.Columns.Item(3).ColumnWidth = 12.5
